$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18-21, columns B-F get their values replaced with "-" to remove
# duplicated professor lists.
$ws.Range("B18:F21").Value = "-"
